# Update workbook "江西-漫展信息.xlsx" to the newer scraped data snapshot.
#
# Net effect per sheet:
#  - 展览 (Exhibitions) and 全部类型 (All types): the two oldest events
#    (上饶·囧喵喵国风动漫展 id=88050 and 南昌·第一届异次元动漫嘉年华 id=84102)
#    have fallen out of the listing window and are removed; every remaining
#    row shifts up by two, the serial numbers in column A are renumbered,
#    and the "want to go" counters (column F) are refreshed for several
#    still-listed events.
#  - 演出 (Shows): the "want to go" counter for the CrossingX show is
#    refreshed as well.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------------
# 1. Remove the two obsolete rows (rows 2 and 3) from 展览 and 全部类型.
#    EntireRow.Delete() shifts everything below up by two rows, and Excel
#    recalculates each sheet's used-range / <dimension> automatically.
# ---------------------------------------------------------------------------
$ws1.Range("A2:A3").EntireRow.Delete()
$ws4.Range("A2:A3").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Renumber column A (the 0-based serial number column) for both sheets
#    so it again reads 0,1,2,3,... with no gaps.
#    (UsedRange is used instead of CurrentRegion, which does not reliably
#    recompute immediately after a row deletion in this environment.)
# ---------------------------------------------------------------------------
$lastRow1 = $ws1.UsedRange.Rows.Count
for ($i = 2; $i -le $lastRow1; $i++) {
    $ws1.Cells.Item($i, 1).Value = $i - 1
}

$lastRow4 = $ws4.UsedRange.Rows.Count
for ($i = 2; $i -le $lastRow4; $i++) {
    $ws4.Cells.Item($i, 1).Value = $i - 1
}

# ---------------------------------------------------------------------------
# 3. Refresh the "want to go" (想去人数, column F) counts that changed for
#    the events that are still listed, on both 展览 and 全部类型.
# ---------------------------------------------------------------------------
$updatesF1 = @{
    3  = 143
    4  = 1743
    8  = 31
    9  = 61
    13 = 143
    18 = 4599
    19 = 42
    20 = 815
    22 = 2179
    23 = 70
    24 = 12
    25 = 2040
}
foreach ($row in $updatesF1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updatesF1[$row]
}

$updatesF4 = @{
    3  = 143
    4  = 1743
    8  = 31
    9  = 61
    13 = 143
    18 = 4599
    19 = 73
    20 = 42
    22 = 815
    24 = 2179
    25 = 70
    26 = 12
    27 = 2040
}
foreach ($row in $updatesF4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updatesF4[$row]
}

# ---------------------------------------------------------------------------
# 4. 演出 sheet: refresh the "want to go" count for the CrossingX show.
# ---------------------------------------------------------------------------
$ws2.Cells.Item(2, 6).Value = 73
